$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap values of J3 and J4
$j3 = $ws.Range("J3").Value2
$j4 = $ws.Range("J4").Value2
$ws.Range("J3").Value = $j4
$ws.Range("J4").Value = $j3

# Update selection to K5
$ws.Range("K5").Select()
